$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.016.46'
$ws.Range("E2").Value = '  -1.91%  '
$ws.Range("D3").Value = '3.726.30'
$ws.Range("E3").Value = '  -2.27%  '
$ws.Range("E4").Value = '  +0.22%  '
$ws.Range("D5").Value = "'620.44"
$ws.Range("E5").Value = '  -0.10%  '
$ws.Range("E6").Value = '  -0.57%  '
$ws.Range("D7").Value = '3.725.91'
$ws.Range("E7").Value = '  -2.03%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("E9").Value = '  -1.92%  '
$ws.Range("D10").Value = "'0.166"
$ws.Range("E10").Value = '  +0.58%  '
$ws.Range("E11").Value = '  -4.30%  '
$ws.Range("D12").Value = "'0.484"
$ws.Range("E12").Value = '  -3.81%  '
$ws.Range("D13").Value = "'40.53"
$ws.Range("E13").Value = '  -0.36%  '
$ws.Range("D14").Value = "'0.0000257"
$ws.Range("E14").Value = '  -0.10%  '
$ws.Range("D15").Value = '4.351.32'
$ws.Range("E15").Value = '  -1.55%  '
$ws.Range("D16").Value = '3.733.58'
$ws.Range("E16").Value = '  -1.62%  '
$ws.Range("D17").Value = '70.034.29'
$ws.Range("E17").Value = '  -1.91%  '
$ws.Range("E18").Value = '  -1.86%  '
$ws.Range("E19").Value = '  +0.23%  '
$ws.Range("D20").Value = "'16.67"
$ws.Range("E20").Value = '  -1.99%  '
$ws.Range("D21").Value = "'505.03"
$ws.Range("E21").Value = '  -3.34%  '
$ws.Range("D22").Value = "'9.28"
$ws.Range("E22").Value = '  -1.72%  '
$ws.Range("E23").Value = '  -4.55%  '
$ws.Range("D24").Value = "'2.54"
$ws.Range("E24").Value = '  -0.15%  '
$ws.Range("D25").Value = "'86.52"
$ws.Range("E25").Value = '  -3.05%  '
$ws.Range("D26").Value = "'11.43"
$ws.Range("E26").Value = '  +2.05%  '
$ws.Range("D27").Value = "'13.07"
$ws.Range("E27").Value = '  -4.19%  '
$ws.Range("E28").Value = '  +20.01%  '
$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = '  -0.18%  '
$ws.Range("E30").Value = '  -3.33%  '
$ws.Range("E31").Value = '  +0.50%  '
$ws.Range("E32").Value = '  -2.86%  '
$ws.Range("D33").Value = "'31.10"
$ws.Range("E33").Value = '  -4.22%  '
$ws.Range("E34").Value = '  -1.85%  '
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = '  +0.16%  '
$ws.Range("E36").Value = '  -0.31%  '
$ws.Range("D37").Value = "'6.15"
$ws.Range("E37").Value = '  -0.50%  '
$ws.Range("E38").Value = '  +1.92%  '
$ws.Range("D39").Value = "'0.338"
$ws.Range("E39").Value = '  -1.88%  '
$ws.Range("E40").Value = '  -7.65%  '
$ws.Range("D41").Value = "'50.18"
$ws.Range("E41").Value = '  -3.16%  '
$ws.Range("D42").Value = "'45.41"
$ws.Range("E42").Value = '  +1.07%  '
$ws.Range("D43").Value = "'431.82"
$ws.Range("E43").Value = '  -2.83%  '
$ws.Range("D44").Value = "'2.88"
$ws.Range("E44").Value = '  +1.34%  '
$ws.Range("D45").Value = "'8.67"
$ws.Range("E45").Value = '  -2.52%  '
$ws.Range("D46").Value = '2.994.70'
$ws.Range("E46").Value = '  -6.14%  '
$ws.Range("D47").Value = "'0.0362"
$ws.Range("E47").Value = '  -1.55%  '
$ws.Range("D48").Value = "'27.50"
$ws.Range("E48").Value = '  -2.18%  '
$ws.Range("E49").Value = '  +0.00%  '
$ws.Range("D50").Value = "'137.04"
$ws.Range("E50").Value = '  -2.44%  '
$ws.Range("D51").Value = "'2.49"
$ws.Range("E51").Value = '  +0.81%  '
